# Update rebound pollutant emissions for LULUCF
#
# 1) "About" sheet: add a new "CO2 Sequestration Start Year" source block
#    (U.S. State Department / Second Biennial Report) above the existing
#    "Rebound CH4 and N2O Emissions" block, pushing everything below it
#    down by 7 rows.
# 2) "data from RPEpUACE" sheet: update the CH4 and N2O rebound emission
#    factors. Downstream formulas on "BLAPE" recompute automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "About" sheet
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Make room for the new source block: insert 7 blank rows starting at
# row 10 (the old "Rebound CH4 and N2O Emissions" header and everything
# under it shifts from rows 10-22 down to rows 17-29).
$wsAbout.Rows.Item(10).Resize(7).Insert() | Out-Null

# Pull formatting for the new header/year/url cells from the matching
# cells in the (now shifted) "US EPA" source block above, then stamp in
# the new text/values.
$wsAbout.Range("B17").Copy() | Out-Null
$wsAbout.Range("B10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsAbout.Range("B5").Copy() | Out-Null
$wsAbout.Range("B12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsAbout.Range("B7").Copy() | Out-Null
$wsAbout.Range("B14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsAbout.Application.CutCopyMode = 0

$wsAbout.Range("B10").Value = "CO2 Sequestration"
$wsAbout.Range("B11").Value = "U.S. State Department"
$wsAbout.Range("B12").Value = 2016
$wsAbout.Range("B13").Value = "Second Biennial Report of the United States of America"
$wsAbout.Range("B14").Value = "https://unfccc.int/files/national_reports/biennial_reports_and_iar/submitted_biennial_reports/application/pdf/2016_second_biennial_report_of_the_united_states_.pdf"
$wsAbout.Range("B15").Value = "Page 34, Table 3"

# ---------------------------------------------------------------------
# 2) "data from RPEpUACE" sheet - updated rebound emission factors
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data from RPEpUACE")
$wsData.Range("B11").Value = 0.00066417418588509813
$wsData.Range("B12").Value = 0.000038529345186784264

$wsData.Range("B2:B13").Select() | Out-Null

# Leave "About" as the active sheet/tab, matching the saved workbook view.
$wsAbout.Activate() | Out-Null
$wsAbout.Rows.Item(16).Select() | Out-Null
